# Auto-generated: update crypto price (D) and 1h volume change (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'64.862.72"
$ws.Range("E2").Value = "  +2.22%  "
$ws.Range("D3").Value = "'2.640.76"
$ws.Range("E3").Value = "  +2.33%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'594.13"
$ws.Range("E5").Value = "  +0.86%  "
$ws.Range("D6").Value = "'155.33"
$ws.Range("E6").Value = "  +3.52%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").Value = "'0.590"
$ws.Range("E8").Value = "  +0.56%  "
$ws.Range("E9").Value = "  +5.83%  "
$ws.Range("D10").Value = "'0.398"
$ws.Range("E10").Value = "  +3.77%  "
$ws.Range("D11").Value = "'5.78"
$ws.Range("E11").Value = "  +0.82%  "
$ws.Range("E12").Value = "  +1.91%  "
$ws.Range("D13").Value = "'29.04"
$ws.Range("E13").Value = "  +5.24%  "
$ws.Range("D14").Value = "'0.0000186"
$ws.Range("E14").Value = "  +19.72%  "
$ws.Range("D15").Value = "'3.114.22"
$ws.Range("E15").Value = "  +2.29%  "
$ws.Range("D16").Value = "'64.786.15"
$ws.Range("E16").Value = "  +2.34%  "
$ws.Range("D17").Value = "'2.630.78"
$ws.Range("E17").Value = "  +2.15%  "
$ws.Range("D18").Value = "'12.53"
$ws.Range("E18").Value = "  +2.51%  "
$ws.Range("D19").Value = "'4.79"
$ws.Range("E19").Value = "  +1.70%  "
$ws.Range("D20").Value = "'350.73"
$ws.Range("E20").Value = "  +1.20%  "
$ws.Range("D21").Value = "'7.28"
$ws.Range("E21").Value = "  +6.15%  "
$ws.Range("E22").Value = "  +0.16%  "
$ws.Range("D23").Value = "'67.79"
$ws.Range("E23").Value = "  +0.74%  "
$ws.Range("D24").Value = "'1.69"
$ws.Range("E24").Value = "  -0.47%  "
$ws.Range("D25").Value = "'9.48"
$ws.Range("E25").Value = "  +3.62%  "
$ws.Range("D26").Value = "'1.63"
$ws.Range("E26").Value = "  -2.56%  "
$ws.Range("E27").Value = "  +0.10%  "
$ws.Range("E28").Value = "  +0.51%  "
$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "  +0.09%  "
$ws.Range("D30").Value = "'0.0₃0934"
$ws.Range("E30").Value = "  +8.66%  "
$ws.Range("D31").Value = "'2.09"
$ws.Range("E31").Value = "  +2.67%  "
$ws.Range("D32").Value = "'511.08"
$ws.Range("E32").Value = "  -7.90%  "
$ws.Range("E33").Value = "  +0.46%  "
$ws.Range("E34").Value = "  +6.60%  "
$ws.Range("D35").Value = "'6.26"
$ws.Range("E35").Value = "  +3.31%  "
$ws.Range("D36").Value = "'0.424"
$ws.Range("E36").Value = "  +2.55%  "
$ws.Range("D37").Value = "'164.77"
$ws.Range("E37").Value = "  -1.08%  "
$ws.Range("D38").Value = "'20.16"
$ws.Range("E38").Value = "  +3.23%  "
$ws.Range("E39").Value = "  +3.96%  "
$ws.Range("E40").Value = "  +0.02%  "
$ws.Range("E41").Value = "  +0.06%  "
$ws.Range("D42").Value = "'42.19"
$ws.Range("E42").Value = "  +6.11%  "
$ws.Range("D43").Value = "'164.22"
$ws.Range("E43").Value = "  -1.04%  "
$ws.Range("D44").Value = "'4.08"
$ws.Range("E44").Value = "  +1.71%  "
$ws.Range("D45").Value = "'0.0612"
$ws.Range("E45").Value = "  +3.89%  "
$ws.Range("D46").Value = "'22.79"
$ws.Range("E46").Value = "  -1.53%  "
$ws.Range("E47").Value = "  +3.14%  "
$ws.Range("D48").Value = "'0.646"
$ws.Range("E48").Value = "  +2.69%  "
$ws.Range("D49").Value = "'0.0253"
$ws.Range("E49").Value = "  +0.88%  "
$ws.Range("D50").Value = "'0.0979"
$ws.Range("E50").Value = "  +1.65%  "
$ws.Range("E51").Value = "  +0.39%  "
